$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.062.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.018.51'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.89%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '557.60'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.21'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.563'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.025.97'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.365'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.551.32'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.99%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.096.47'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.41%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.022.07'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.89%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000150'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '396.45'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.12'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.88'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.65'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.72%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.21'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.467'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.33%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0₃0978'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.99%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.995'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.48'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '159.75'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.71'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.06'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.76%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.10'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.78%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.521.02'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -6.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.59'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.72'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.53%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.64'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.669'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0600'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.18%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.08'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -6.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0247'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.93'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.78%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0947'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '10.51'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '264.17'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.92%  '
